$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fix the Quebec row's DirectStyleURL (C7): the old Mapbox style id was
#    replaced with a new one.
# ---------------------------------------------------------------------------
$ws.Range("C7").Value = "ckv8617pd0n9g14rxpi45ts4u"

# ---------------------------------------------------------------------------
# 2. Add three new metric columns (AH, AI, AJ): YearOfStats, SourceGTFS,
#    DateUpdatedGTFS. Header cells span rows 1:2 like the other headers.
# ---------------------------------------------------------------------------
$ws.Range("AH1:AH2").Merge()
$ws.Range("AI1:AI2").Merge()
$ws.Range("AJ1:AJ2").Merge()

$ws.Range("AH1").Value = "YearOfStats"
$ws.Range("AI1").Value = "SourceGTFS"
$ws.Range("AJ1").Value = "DateUpdatedGTFS"

$ws.Range("AH1:AJ2").HorizontalAlignment = -4108
$ws.Range("AH1:AJ2").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 3. Fill in the new metric values for every city row.
# ---------------------------------------------------------------------------
# Montreal (row 3)
$ws.Cells.Item(3, 34).Value = 2016
$ws.Cells.Item(3, 35).Value = "https://www.stm.info/en/about/developers"
$ws.Cells.Item(3, 36).Value = 2021

# Vienna (row 4)
$ws.Cells.Item(4, 34).Value = 2019
$ws.Cells.Item(4, 35).Value = "https://www.data.gv.at/katalog/dataset/wiener-linien-fahrplandaten-gtfs-wien"
$ws.Cells.Item(4, 36).Value = 2021

# Barcelona (row 5)
$ws.Cells.Item(5, 34).Value = 2015
$ws.Cells.Item(5, 35).Value = "https://developer.tmb.cat/data"
$ws.Cells.Item(5, 36).Value = 2021

# Budapest (row 6)
$ws.Cells.Item(6, 34).Value = 2014
$ws.Cells.Item(6, 35).Value = "https://data.europa.eu/it/news/centre-budapest-transport-and-open-mobility-data"
$ws.Cells.Item(6, 36).Value = 2021

# Quebec (row 7)
$ws.Cells.Item(7, 34).Value = 2016
$ws.Cells.Item(7, 35).Value = "https://www.rtcquebec.ca/donnees-ouvertes"
$ws.Cells.Item(7, 36).Value = 2021

# Center-align the new data cells, matching the rest of the table.
$ws.Range("AH3:AJ7").HorizontalAlignment = -4108
$ws.Range("AH3:AJ7").VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# 4. Column widths for the new / newly-sized columns.
# ---------------------------------------------------------------------------
$ws.Columns.Item(30).ColumnWidth = 12.66   # AD
$ws.Columns.Item(32).ColumnWidth = 13.44   # AF
$ws.Columns.Item(33).ColumnWidth = 12.55   # AG
$ws.Columns.Item(34).ColumnWidth = 12.22   # AH
$ws.Columns.Item(35).ColumnWidth = 70.66   # AI
$ws.Columns.Item(36).ColumnWidth = 16      # AJ

# ---------------------------------------------------------------------------
# 5. Update the sheet view: scroll position and active selection.
# ---------------------------------------------------------------------------
$ws.Range("AH10").Select()
$excel.ActiveWindow.ScrollColumn = 24
